# Normalize the "Recorded By" (column G) values so that entries of the
# form "System, <name>" become "<name>, System" (moving "System" to the
# end of the comma-separated list). Cells with a different number of
# entries (e.g. single names, or three-part combinations) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $parts = $val.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($first -eq "System") {
                $cell.Value = "$second, $first"
            }
        }
    }
}
